$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.607.40'
$ws.Range('E2').Value = '  +3.22%  '
$ws.Range('D3').Value = '1.797.60'
$ws.Range('E3').Value = '  +0.12%  '
$cell = $ws.Range('D4')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$ws.Range('E4').Value = '  +0.37%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '313.36'
$cell.ClearFormats()
$ws.Range('E5').Value = '  -0.15%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$ws.Range('E6').Value = '  +0.27%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '0.5362'
$cell.ClearFormats()
$ws.Range('E7').Value = '  -1.01%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.3777'
$cell.ClearFormats()
$ws.Range('E8').Value = '  +0.10%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.07534'
$cell.ClearFormats()
$ws.Range('E9').Value = '  -0.03%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '42.51'
$cell.ClearFormats()
$ws.Range('E10').Value = '  -0.58%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '1.116'
$cell.ClearFormats()
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('B12').Value = 'BinanceUSD'
$ws.Range('C12').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '20.97'
$cell.ClearFormats()
$ws.Range('E13').Value = '  +0.13%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '6.188'
$cell.ClearFormats()
$ws.Range('E14').Value = '  -0.27%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '7.446'
$cell.ClearFormats()
$ws.Range('E15').Value = '  +4.98%  '
$ws.Range('D16').Value = '1.797.55'
$ws.Range('E16').Value = '  +0.47%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '90.29'
$cell.ClearFormats()
$ws.Range('E17').Value = '  -0.46%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '0.00001065'
$cell.ClearFormats()
$ws.Range('E18').Value = '  -0.56%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '0.06448'
$cell.ClearFormats()
$ws.Range('E19').Value = '  -0.82%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.ClearFormats()
$ws.Range('E20').Value = '  +0.18%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '17.23'
$cell.ClearFormats()
$ws.Range('E21').Value = '  +1.27%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '5.930'
$cell.ClearFormats()
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D23').Value = '28.616.82'
$ws.Range('E23').Value = '  +3.08%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '11.17'
$cell.ClearFormats()
$ws.Range('E24').Value = '  -0.48%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '2.098'
$cell.ClearFormats()
$ws.Range('E25').Value = '  -0.29%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '160.63'
$cell.ClearFormats()
$ws.Range('E26').Value = '  +3.67%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '20.46'
$cell.ClearFormats()
$ws.Range('E27').Value = '  -0.10%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '2.382'
$cell.ClearFormats()
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('D29').Value = '2.003.76'
$ws.Range('E29').Value = '  +0.36%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '123.05'
$cell.ClearFormats()
$ws.Range('E30').Value = '  +0.45%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '1.107'
$cell.ClearFormats()
$ws.Range('E31').Value = '  -3.54%  '
$ws.Range('E32').Value = '  -0.58%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '5.668'
$cell.ClearFormats()
$ws.Range('E33').Value = '  -0.40%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '3.682'
$cell.ClearFormats()
$ws.Range('E34').Value = '  +1.81%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '0.2254'
$cell.ClearFormats()
$ws.Range('E35').Value = '  +7.42%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '0.06463'
$cell.ClearFormats()
$ws.Range('E36').Value = '  +7.17%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '8.882'
$cell.ClearFormats()
$ws.Range('E37').Value = '  +2.85%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '0.02300'
$cell.ClearFormats()
$ws.Range('E38').Value = '  +0.61%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '5.042'
$cell.ClearFormats()
$ws.Range('E39').Value = '  +0.73%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '11.36'
$cell.ClearFormats()
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '1.211'
$cell.ClearFormats()
$ws.Range('E41').Value = '  +5.52%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '0.6255'
$cell.ClearFormats()
$ws.Range('E42').Value = '  -0.03%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.ClearFormats()
$ws.Range('E43').Value = '  +0.12%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '1.391'
$cell.ClearFormats()
$ws.Range('E44').Value = '  -0.55%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '13.48'
$cell.ClearFormats()
$ws.Range('E45').Value = '  +0.89%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '0.5885'
$cell.ClearFormats()
$ws.Range('E46').Value = '  +0.20%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '3.661'
$cell.ClearFormats()
$ws.Range('E47').Value = '  +0.91%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '125.89'
$cell.ClearFormats()
$ws.Range('E48').Value = '  +3.31%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '1.964'
$cell.ClearFormats()
$ws.Range('E49').Value = '  +2.31%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '1.158'
$cell.ClearFormats()
$ws.Range('E50').Value = '  +2.19%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '0.06909'
$cell.ClearFormats()
$ws.Range('E51').Value = '  +1.92%  '
